# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 4f2cba96-... file row on the zh-cn and de-de report
# sheets, reflecting a newer handback run.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-18 02:49:10"
$wsZhCn.Range("H3").Value = "2016-03-18 02:49:55"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-18 02:49:18"
$wsDeDe.Range("H3").Value = "2016-03-18 02:50:10"
